# Daily attendance processing - 2025-10-25 23:41:30
# Swap the order of the two comma-separated "Recorded By" entries in column G
# for the specific rows where the recorded-by list has exactly two names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,6,7,10,11,12,13,14,15,17,18,19,20,21,22,30,33,34,37,38,39,40,41,42,44,45,46,47,48,49,57,60,61,64,65,66,67,68,69,71,72,73,74,75,76,86,87,88,89,90,93,95,96,97,99,102,112,113,114,115,116,119,121,122,123,125,128,138,139,140,141,142,145,147,148,149,151,154)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value
    $parts = $current -split ", ", 2
    if ($parts.Count -eq 2) {
        $cell.Value = $parts[1] + ", " + $parts[0]
    }
}
